$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.451.85"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.574.11"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "'287.99"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.3727"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("D8").Value = "'47.53"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'0.3322"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'1.157"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("D11").Value = "'0.07519"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "1.566.31"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'0.00001119"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "'88.36"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'0.06730"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'6.396"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'16.52"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "22.444.49"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'2.393"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "'2.626"
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("D27").Value = "'150.60"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'19.64"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "'4.955"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "'125.29"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "1.744.47"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'1.099"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "'6.099"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "'1.989"
$ws.Range("D35").Value = "'9.843"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").Value = "'0.08338"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").Value = "'0.02460"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("D38").Value = "'1.316"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").Value = "'0.06398"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "'5.350"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "'0.6272"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.07"
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'0.6094"
$ws.Range("D47").Value = "'3.776"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'2.049"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "'125.03"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  -0.14%  "
